$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "'1744624259342"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "vishal K"
$ws.Range("C5").Value = "vishal@ssn.edu"
$ws.Range("D5").Value = '$2a$10$WFBDROOHGnBjzR7XpSrwrOQjI8BjASxXN67YRbfYB2sDga9YyAlO6'
$ws.Range("E5").Value = "CSE"
$ws.Range("F5").Value = "2025-04-14T09:50:59.342Z"

# Row 6
$ws.Range("A6").Value = "'1744625892798"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "test user"
$ws.Range("C6").Value = "test1@college.edu"
$ws.Range("D6").Value = '$2a$10$GRSNPeX9QoxldqbiG5BRSOxfRQVFVJ42AUrpnUA1MUhDfvOQjuH3.'
$ws.Range("E6").Value = "CSE"
$ws.Range("F6").Value = "2025-04-14T10:18:12.798Z"
